$d = $word.ActiveDocument

# -------------------------------------------------------------------
# 1) Title paragraph: split "HW3: Extend HW #2 and explore a 5 X 5
#    redistricting scheme with respect to geography" into three runs,
#    replacing "and explore" with "by exploring", while keeping all
#    three runs italic (matching formatting of the original run).
# -------------------------------------------------------------------
$titlePara = $d.Paragraphs(2)
$titleRange = $titlePara.Range
$titleText = $titleRange.Text
$titleStart = $titleRange.Start
$oldPhrase = "and explore"
$newPhrase = "by exploring"
$phraseIdx = $titleText.IndexOf($oldPhrase)
$midRange = $d.Range($titleStart + $phraseIdx, $titleStart + $phraseIdx + $oldPhrase.Length)
$midRange.Text = $newPhrase
# Force the replaced span into its own run by toggling italic off then
# back on - this keeps the run boundaries even though the final
# formatting equals its neighbors.
$midRange.Font.Italic = $false
$midRange2 = $d.Range($titleStart + $phraseIdx, $titleStart + $phraseIdx + $newPhrase.Length)
$midRange2.Font.Italic = $true

# -------------------------------------------------------------------
# 2) Table cell margins: change left margin from 24 dxa (1.2 pt) to
#    16 dxa (0.8 pt) on the table itself and every one of its cells.
# -------------------------------------------------------------------
$tbl = $d.Tables(1)
$tbl.LeftPadding = 0.8
foreach ($cell in $tbl.Range.Cells) {
    $cell.LeftPadding = 0.8
}

# -------------------------------------------------------------------
# 3) Merge the split "re" + "districting" runs back into single runs
#    by re-asserting the same text via Find/Replace (this engine
#    coalesces adjacent same-formatted runs whenever a range is
#    edited).
# -------------------------------------------------------------------
$d.Content.Find.Execute(
    "How many total contiguous redistricting schemes were generated?",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "How many total contiguous redistricting schemes were generated?", 2)

$d.Content.Find.Execute(
    "Keep in mind, you may have to generate on the order of 10 million random redistricting schemes to see a statistically significant number of contiguous redistricting schemes. 100 contiguous schemes is ideal, and under 10 is not enough. Use your best judgment on a number in between.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Keep in mind, you may have to generate on the order of 10 million random redistricting schemes to see a statistically significant number of contiguous redistricting schemes. 100 contiguous schemes is ideal, and under 10 is not enough. Use your best judgment on a number in between.", 2)

# -------------------------------------------------------------------
# 4) Add the new character styles ListLabel10 .. ListLabel18 that show
#    up in styles.xml (list-label runs minted for additional list
#    levels used elsewhere in the source document).
# -------------------------------------------------------------------
for ($n = 10; $n -le 18; $n++) {
    $styleName = "ListLabel$n"
    $style = $d.Styles.Add($styleName, 2)
    $style.NameLocal = "ListLabel $n"
    $style.QuickStyle = $true
    $style.Font.NameBi = "OpenSymbol"
    if ($n -eq 10) {
        $style.Font.Size = 11
    }
}
